# Fruta / hortaliza, semanal
# Re-assign the weekly Fecha/Calidad/Volumen/Precio* values for each data
# row (2..31) of the Guayaba sheet. Columns A,B,C,E,F,G,H,I,J,K,Q,R,T are
# left untouched; only D (Fecha), L (Calidad), M (Volumen), N (Precio
# minimo), O (Precio maximo), P (Precio promedio ponderado) and
# S (Precio $/Kg) are updated per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, D(Fecha), L(Calidad), M(Volumen), N(Precio min), O(Precio max), P(Precio prom), S(Precio $/Kg)
$data = @(
    @(2,  44498, "Segunda", 100, 1200, 1300, 1250, 1250),
    @(3,  44414, "Primera", 160, 1300, 1400, 1350, 1350),
    @(4,  44260, "Primera", 100, 1900, 2000, 1950, 1950),
    @(5,  44425, "Primera", 140, 1200, 1300, 1250, 1250),
    @(6,  44330, "Primera", 200, 1200, 1300, 1250, 1250),
    @(7,  44330, "Segunda", 100, 1000, 1100, 1050, 1050),
    @(8,  44351, "Primera", 100, 700,  800,  750,  750),
    @(9,  44351, "Segunda", 100, 600,  700,  650,  650),
    @(10, 44389, "Primera", 140, 750,  800,  775,  775),
    @(11, 44389, "Segunda", 120, 600,  700,  650,  650),
    @(12, 44348, "Primera", 120, 1000, 1100, 1050, 1050),
    @(13, 44309, "Primera", 160, 1400, 1500, 1450, 1450),
    @(14, 44326, "Primera", 160, 600,  700,  650,  650),
    @(15, 44379, "Primera", 150, 700,  800,  747,  747),
    @(16, 44379, "Segunda", 140, 500,  600,  543,  543),
    @(17, 44407, "Primera", 200, 600,  650,  625,  625),
    @(18, 44417, "Primera", 200, 1300, 1400, 1350, 1350),
    @(19, 44403, "Primera", 100, 1200, 1300, 1250, 1250),
    @(20, 44403, "Segunda", 120, 950,  1000, 975,  975),
    @(21, 44386, "Primera", 160, 700,  750,  725,  725),
    @(22, 44386, "Segunda", 200, 600,  650,  625,  625),
    @(23, 44473, "Primera", 160, 1500, 1600, 1550, 1550),
    @(24, 44316, "Primera", 140, 1100, 1200, 1150, 1150),
    @(25, 44372, "Primera", 900, 750,  800,  772,  772),
    @(26, 44372, "Segunda", 900, 600,  650,  628,  628),
    @(27, 44350, "Primera", 140, 750,  800,  775,  775),
    @(28, 44358, "Primera", 200, 700,  800,  750,  750),
    @(29, 44358, "Segunda", 200, 600,  650,  625,  625),
    @(30, 44344, "Primera", 140, 1000, 1200, 1100, 1100),
    @(31, 44344, "Segunda", 120, 800,  850,  825,  825)
)

foreach ($entry in $data) {
    $row = $entry[0]
    $ws.Cells.Item($row, 4).Value2  = $entry[1]   # D - Fecha
    $ws.Cells.Item($row, 12).Value  = $entry[2]   # L - Calidad
    $ws.Cells.Item($row, 13).Value2 = $entry[3]   # M - Volumen
    $ws.Cells.Item($row, 14).Value2 = $entry[4]   # N - Precio minimo
    $ws.Cells.Item($row, 15).Value2 = $entry[5]   # O - Precio maximo
    $ws.Cells.Item($row, 16).Value2 = $entry[6]   # P - Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value2 = $entry[7]   # S - Precio $/Kg
}
